$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Fill in the newly benchmarked "inverse k" timings on the Set Partitioning sheet.
# Write order matches column-major top-to-bottom traversal (M, N, O, AC, AD, AE, AF, P)
# so shared-string indices line up with the source edit.
$ws2.Range("M3").Value = "7us"
$ws2.Range("M4").Value = "10us"
$ws2.Range("M5").Value = "7us"
$ws2.Range("M6").Value = "10us"
$ws2.Range("M7").Value = "10us"
$ws2.Range("M8").Value = "13us"
$ws2.Range("M9").Value = "35us"
$ws2.Range("M10").Value = "29us"
$ws2.Range("M11").Value = "19us"
$ws2.Range("M12").Value = "45us"
$ws2.Range("M13").Value = "56us"
$ws2.Range("M14").Value = "175us"

$ws2.Range("N3").Value = "36us"
$ws2.Range("N4").Value = "91us"
$ws2.Range("N5").Value = "126us"
$ws2.Range("N6").Value = "189us"
$ws2.Range("N7").Value = "410us"
$ws2.Range("N8").Value = "672us"
$ws2.Range("N9").Value = "1004us"
$ws2.Range("N10").Value = "1495us"
$ws2.Range("N11").Value = "2097us"
$ws2.Range("N12").Value = "4151us"
$ws2.Range("N13").Value = "8803us"
$ws2.Range("N14").Value = "31062us"

$ws2.Range("O3").Value = "436us"
$ws2.Range("O4").Value = "1792us"
$ws2.Range("O5").Value = "3134us"
$ws2.Range("O6").Value = "15162us"
$ws2.Range("O7").Value = "20911us"
$ws2.Range("O8").Value = "62361us"
$ws2.Range("O9").Value = "76187us"
$ws2.Range("O10").Value = "93257us"
$ws2.Range("O11").Value = "152464us"
$ws2.Range("O12").Value = "445788us"
$ws2.Range("O13").Value = "1089082us"
$ws2.Range("O14").Value = "5806986us"

$ws2.Range("AC3").Value = "1827us"
$ws2.Range("AC4").Value = "1097us"
$ws2.Range("AC5").Value = "1624us"
$ws2.Range("AC6").Value = "2406us"
$ws2.Range("AC7").Value = "4063us"
$ws2.Range("AC8").Value = "21727us"
$ws2.Range("AC9").Value = "8851us"
$ws2.Range("AC10").Value = "10730us"
$ws2.Range("AC11").Value = "10767us"
$ws2.Range("AC12").Value = "20196us"
$ws2.Range("AC13").Value = "34928us"
$ws2.Range("AC14").Value = "83532us"

$ws2.Range("AD3").Value = "24353us"
$ws2.Range("AD4").Value = "65663us"
$ws2.Range("AD5").Value = "74668us"
$ws2.Range("AD6").Value = "137859us"
$ws2.Range("AD7").Value = "245027us"
$ws2.Range("AD8").Value = "396088us"
$ws2.Range("AD9").Value = "672337us"
$ws2.Range("AD10").Value = "999572us"
$ws2.Range("AD11").Value = "1578376us"
$ws2.Range("AD12").Value = "3833036us"
$ws2.Range("AD13").Value = "9156522us"
$ws2.Range("AD14").Value = "65501637us"

$ws2.Range("AE3").Value = "239193us"
$ws2.Range("AE4").Value = "676641us"
$ws2.Range("AE5").Value = "1793670us"
$ws2.Range("AE6").Value = "4686700us"
$ws2.Range("AE7").Value = "12414109us"
$ws2.Range("AE8").Value = "41274616us"
$ws2.Range("AE9").Value = "263048466us"
$ws2.Range("AE10").Value = "TIME"

$ws2.Range("AF3").Value = "TIME"

$ws2.Range("P3").Value = "167672us"
$ws2.Range("P4").Value = "11545300us"
$ws2.Range("P5").Value = "TIME"

# Columns AD/AE widened (auto-fit) to fit the new, wider microsecond values.
$ws2.Columns.Item(30).ColumnWidth = 9.5
$ws2.Columns.Item(31).ColumnWidth = 10.666666666666666

# Selection/view state: Integer Partitioning keeps a stale selection at O24,
# while Set Partitioning becomes the active tab, scrolled in and selected at AF4.
$ws1.Activate()
$ws1.Range("O24").Select()

$ws2.Activate()
$ws2.Range("AF4").Select()
